$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Select()
$ws.Columns("N").Insert()
$ws.Columns("N").Select()
